# Atualização de bases das ligas, do dia: 28-04-2024 às 15:37
#
# Sheet "Mexico Liga de Expansion" updates:
#  1) Rows 186/187 (match ids 7648957 / 7648958) had their results/odds
#     swapped (including which teams go with which stats).
#  2) Row 239 (match id 8127905 -> 8127903) moved one day later and its
#     teams / odds were corrected.
#  3) Row 240 (match id 8127904 -> 8127899) moved one day later and its
#     teams / odds were corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 186 / 187 : swap everything except the running index (A),
#     the league (C) and the kickoff date (D). -----------------------

# id (B)
$b186 = $ws.Range("B186").Value()
$b187 = $ws.Range("B187").Value()
$ws.Range("B186").Value = $b187
$ws.Range("B187").Value = $b186

# HomeTeam / AwayTeam (E / F)
$e186 = $ws.Range("E186").Value()
$f186 = $ws.Range("F186").Value()
$e187 = $ws.Range("E187").Value()
$f187 = $ws.Range("F187").Value()
$ws.Range("E186").Value = $e187
$ws.Range("F186").Value = $f187
$ws.Range("E187").Value = $e186
$ws.Range("F187").Value = $f186

# FTHG..PL_AhUnder (G:AB) - odds & results
$g186ab = $ws.Range("G186:AB186").Value()
$g187ab = $ws.Range("G187:AB187").Value()
$ws.Range("G186:AB186").Value = $g187ab
$ws.Range("G187:AB187").Value = $g186ab

# --- Row 239 : id 8127905 -> 8127903 -------------------------------
# id (B239) is stored as text in this workbook - force text storage
# (via a temporary "@" number format) so we don't turn it into a number,
# then clear the format again so no stray style sticks to the cell.
$ws.Range("B239").NumberFormat = "@"
$ws.Range("B239").Value = "8127903"
$ws.Range("B239").ClearFormats()
$ws.Range("D239").Value = 45410.83333333334
$ws.Range("E239").Value = "Atlante"
$ws.Range("F239").Value = "Cancun FC"
$ws.Range("J239").Value = 1.8
$ws.Range("K239").Value = 3.25
$ws.Range("L239").Value = 4.333
$ws.Range("M239").Value = 1.7
$ws.Range("N239").Value = 3.4
$ws.Range("O239").Value = 5.5
$ws.Range("P239").Value = -0.75
$ws.Range("Q239").Value = 1.9
$ws.Range("R239").Value = 1.9
$ws.Range("S239").Value = 2
$ws.Range("T239").Value = 1.9
$ws.Range("U239").Value = 1.9

# --- Row 240 : id 8127904 -> 8127899 -------------------------------
$ws.Range("B240").NumberFormat = "@"
$ws.Range("B240").Value = "8127899"
$ws.Range("B240").ClearFormats()
$ws.Range("D240").Value = 45410.92013888889
$ws.Range("E240").Value = "Universidad Guadalajara"
$ws.Range("F240").Value = "Mineros de Zacatecas"
$ws.Range("J240").Value = 2.1
$ws.Range("K240").Value = 3.4
$ws.Range("L240").Value = 3.1
$ws.Range("M240").Value = 2.1
$ws.Range("N240").Value = 3.4
$ws.Range("O240").Value = 3.5
$ws.Range("P240").Value = -0.25
$ws.Range("Q240").Value = 1.775
$ws.Range("R240").Value = 2.025
$ws.Range("S240").Value = 2.75
